# "Fruta / hortaliza, semanal"
#
# A new weekly price observation is inserted as row 85 (pushing every
# existing record from row 85 downward by one row, so the former row 227
# becomes row 228). The new row carries the same market/category/quality
# metadata as the record that used to sit at row 85, but with its own
# date and volume.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 85..227 down to 86..228, leaving row 85 empty (formatting,
# e.g. the date style on column D, is carried along by Insert()).
$ws.Rows(85).Insert()

# Populate the newly inserted row 85 with the new observation.
$ws.Cells.Item(85, 1).Value  = 10
$ws.Cells.Item(85, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(85, 3).Value  = "La Araucanía"
$ws.Cells.Item(85, 4).Value  = 44915
$ws.Cells.Item(85, 5).Value  = 9
$ws.Cells.Item(85, 6).Value  = 100114007
$ws.Cells.Item(85, 7).Value  = "Jengibre"
$ws.Cells.Item(85, 8).Value  = "Sin especificar"
$ws.Cells.Item(85, 9).Value  = "Primera"
$ws.Cells.Item(85, 10).Value = 20
$ws.Cells.Item(85, 11).Value = 20000
$ws.Cells.Item(85, 12).Value = 20000
$ws.Cells.Item(85, 13).Value = 20000
$ws.Cells.Item(85, 14).Value = "`$/caja 13 kilos"
$ws.Cells.Item(85, 15).Value = "Perú"
$ws.Cells.Item(85, 16).Value = 1538
$ws.Cells.Item(85, 17).Value = 13
$ws.Cells.Item(85, 18).Value = "Hortaliza"
